$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 44632
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = "Proper Documentation"

$ws.Range("N29").Select() | Out-Null
